$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'22.068.50"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.11%  '

# Row 3
$ws.Range('D3').Value = "'1.554.94"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.16%  '

# Row 4
$ws.Range('D4').Value = "'1.000"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('E5').Value = '  +0.09%  '

# Row 6
$ws.Range('D6').Value = "'291.08"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.33%  '

# Row 7
$ws.Range('D7').Value = "'0.3934"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.40%  '

# Row 8
$ws.Range('D8').Value = "'0.3224"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.81%  '

# Row 9
$ws.Range('D9').Value = "'44.37"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.92%  '

# Row 10
$ws.Range('D10').Value = "'0.07216"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.01%  '

# Row 11
$ws.Range('D11').Value = "'1.076"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.27%  '

# Row 12
$ws.Range('D12').Value = "'1.000"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.06%  '

# Row 13
$ws.Range('D13').Value = "'5.662"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.34%  '

# Row 14
$ws.Range('D14').Value = "'18.70"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.62%  '

# Row 15
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = "'6.694"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.89%  '

# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = "'1.556.12"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.37%  '

# Row 17
$ws.Range('E17').Value = '  +2.31%  '

# Row 18
$ws.Range('D18').Value = "'0.06589"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.53%  '

# Row 19
$ws.Range('D19').Value = "'83.53"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.52%  '

# Row 20
$ws.Range('D20').Value = "'1.000"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.08%  '

# Row 21
$ws.Range('D21').Value = "'6.231"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.38%  '

# Row 22
$ws.Range('D22').Value = "'15.48"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.41%  '

# Row 23
$ws.Range('D23').Value = "'11.22"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.83%  '

# Row 24
$ws.Range('D24').Value = "'22.070.91"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.12%  '

# Row 25
$ws.Range('D25').Value = "'2.357"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.54%  '

# Row 26
$ws.Range('D26').Value = "'2.398"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.84%  '

# Row 27
$ws.Range('D27').Value = "'147.99"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.47%  '

# Row 28
$ws.Range('D28').Value = "'18.53"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.98%  '

# Row 29
$ws.Range('D29').Value = "'4.876"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.08%  '

# Row 30
$ws.Range('D30').Value = "'1.729.23"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.65%  '

# Row 31
$ws.Range('D31').Value = "'118.81"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.09%  '

# Row 32
$ws.Range('D32').Value = "'0.9866"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.62%  '

# Row 33
$ws.Range('D33').Value = "'5.909"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.45%  '

# Row 34
$ws.Range('D34').Value = "'0.08310"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.84%  '

# Row 35
$ws.Range('D35').Value = "'9.146"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.04%  '

# Row 36
$ws.Range('D36').Value = "'1.611"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -15.12%  '

# Row 37
$ws.Range('D37').Value = "'0.02265"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.19%  '

# Row 38
$ws.Range('D38').Value = "'5.123"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.68%  '

# Row 39
$ws.Range('D39').Value = "'0.06026"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.22%  '

# Row 40
$ws.Range('D40').Value = "'1.207"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.92%  '

# Row 41
$ws.Range('D41').Value = "'0.2057"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.34%  '

# Row 42
$ws.Range('D42').Value = "'1.001"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.01%  '

# Row 43
$ws.Range('D43').Value = "'10.73"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.38%  '

# Row 44
$ws.Range('D44').Value = "'0.5806"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.30%  '

# Row 45
$ws.Range('D45').Value = "'3.759"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.75%  '

# Row 46
$ws.Range('E46').Value = '  -4.57%  '

# Row 47
$ws.Range('D47').Value = "'0.5570"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.05%  '

# Row 48
$ws.Range('D48').Value = "'117.62"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.57%  '

# Row 49
$ws.Range('D49').Value = "'1.881"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.50%  '

# Row 50
$ws.Range('D50').Value = "'1.136"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.38%  '

# Row 51
$ws.Range('D51').Value = "'0.06813"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.95%  '
